$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: propagate the "last row" look (bottom border row) that currently
#            sits on row 28/29 onto the new final row (35), before that style
#            gets reassigned to rows 28/29 below.
$ws.Range("C28").Copy()
$ws.Range("C35").PasteSpecial(-4122)

$ws.Range("D28:E28").Copy()
$ws.Range("D35:E35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 2: re-style rows 28 and 29 to match the regular interior row style
#            (same as row 4 and the rest of the table).
$ws.Range("C4").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("D4:E4").Copy()
$ws.Range("D28:E28").PasteSpecial(-4122)
$ws.Range("D29:E29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: fill in the new rows 30-35 (R.B / Faza / Status columns).
# Column D (the shared "Faza" text) is written in the same order the
# original strings were first introduced, so the shared-string table layout
# matches: company-loading, department-loading, department-display,
# company-display, employee-loading, employee-display.
$ws.Range("D30").Value = "Loading sa backend-a za company "
$ws.Range("D32").Value = "Loading sa backend-a za department"
$ws.Range("D33").Value = "Prikaz za department na frontend-u"
$ws.Range("D31").Value = "Prikaz za company na frontend-u"
$ws.Range("D34").Value = "Loading sa backend-a za employee"
$ws.Range("D35").Value = "Prikaz za employees na frontend-u"

$ws.Range("C30").Value = 14.1
$ws.Range("E30").Value = "DONE"

$ws.Range("C31").Value = 14.2
$ws.Range("E31").Value = "DONE"

$ws.Range("C32").Value = 15.1
$ws.Range("E32").Value = "DONE"

$ws.Range("C33").Value = 15.2
$ws.Range("E33").Value = "DONE"

$ws.Range("C34").Value = 16.1
$ws.Range("E34").Value = "DONE"

$ws.Range("C35").Value = 16.2
$ws.Range("E35").Value = "DONE"

# Apply the regular row style (same as row 4) to the new rows 30-34
# (row 35 keeps the "last row" style applied in Step 1).
$ws.Range("C4").Copy()
$ws.Range("C30:C34").PasteSpecial(-4122)

$ws.Range("D4:E4").Copy()
$ws.Range("D30:E34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 4: update the active selection shown in the sheet view
$ws.Range("D32").Select()
